# Append one row (row 15) of trip-log data to the bottom of the sheet,
# matching the diff: dimension A1:H14 -> A1:H15 plus the new row's values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 15

# Column A is blank for this entry, like most rows above it. Assigning ""
# directly would just clear the cell (no stored value at all), so we type
# a lone apostrophe - Excel's "force text" entry - to get a real, empty
# Text cell, then reset the style so we don't leave a stray quote-prefix
# format on the cell (the source rows carry no cell-level style either).
$ws.Cells.Item($row, 1).Value = "'"
$ws.Cells.Item($row, 1).Style = "Normal"

$ws.Cells.Item($row, 2).Value = "احمد"

# Column C is "200", a number-looking value that must stay text (the
# whole column is flagged numberStoredAsText in this workbook). Same
# apostrophe-then-reset-style trick as column A avoids Excel silently
# converting it to a numeric 200.
$ws.Cells.Item($row, 3).Value = "'200"
$ws.Cells.Item($row, 3).Style = "Normal"

$ws.Cells.Item($row, 4).Value = "ايتا"
$ws.Cells.Item($row, 5).Value = "الرحلة 1"
$ws.Cells.Item($row, 6).Value = "C2"
$ws.Cells.Item($row, 7).Value = "NRC"
$ws.Cells.Item($row, 8).Value = "٠٦:٠٣:٤٨ م"
